# Update Fonds de solidarite workbook with 2022-05-31 data.
# Only columns C (nombre_aides) and E (montant_total) change for the
# listed rows; column D (nombre_entreprises) is unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 63;  C = 14355;  E = 36182831 },
    @{ Row = 81;  C = 17433;  E = 133950839 },
    @{ Row = 91;  C = 151134; E = 482406764 },
    @{ Row = 92;  C = 409105; E = 1595107955 },
    @{ Row = 93;  C = 209568; E = 1308770817 },
    @{ Row = 94;  C = 94193;  E = 917466940 },
    @{ Row = 95;  C = 50759;  E = 932280796 },
    @{ Row = 96;  C = 17268;  E = 792060406 },
    @{ Row = 104; C = 135239; E = 272198155 },
    @{ Row = 111; C = 116;    E = 7791001 },
    @{ Row = 115; C = 11695;  E = 32962731 }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 3).Value = $u.C
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
